$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Relay Mosfets" entry (row 12) and the rows that followed it
# (510 Ohm Resistor, 510K Ohm Resistor, 8MHz Clock + its datasheet link).
$ws.Range("A12:D15").ClearContents()

# Re-add "Relay Mosfets" higher up (row 8) with the new N-channel FET part.
$ws.Range("A8").Value = "Relay Mosfets"
$ws.Range("B8").Value = "BUK9Y12-40E"
$ws.Range("D8").Value = "https://assets.nexperia.com/documents/data-sheet/BUK9Y12-40E.pdf"

# New "Necessary Inputs" section listing the required signal connections.
$ws.Range("A18").Value = "Necessary Inputs"
$ws.Range("B18").Value = "DIR"
$ws.Range("A19").Value = "BATT+"
$ws.Range("A20").Value = "BATT GND"
$ws.Range("A21").Value = "EFI_MAIN_SIG"
$ws.Range("A22").Value = "PUMP_SIG"
$ws.Range("A23").Value = "FAN_SIG"
$ws.Range("A24").Value = "INJ_SIG"
$ws.Range("A25").Value = "COIL_SIG"
$ws.Range("A26").Value = "CANH"
$ws.Range("A27").Value = "CANL"
$ws.Range("A28").Value = "AUX"
$ws.Range("A29").Value = "AUX_2"
$ws.Range("A30").Value = "O2"

$ws.Range("A31").Select()
